# Append the new resale-numbers row (2025-01-22 09:02:09) as row 37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

# Columns A & D look like dates/numbers ("2025-01-22", "03") so Excel's
# smart-entry would otherwise coerce them to a date serial / number.
# Force a Text number format for the assignment, then clear the format
# again so the cell ends up with the same (default) style as its
# neighbours, matching how the rest of the sheet is authored.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-22"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "09:02:09"
$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "03"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 126265
$ws.Cells.Item($row, 6).Value = 142142
$ws.Cells.Item($row, 7).Value = 168510
$ws.Cells.Item($row, 8).Value = 158580
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142913
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192230
$ws.Cells.Item($row, 14).Value = 115674
$ws.Cells.Item($row, 15).Value = 45566
$ws.Cells.Item($row, 16).Value = 28458
$ws.Cells.Item($row, 17).Value = 65602
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47931
$ws.Cells.Item($row, 20).Value = -1
